$d = $word.ActiveDocument

# Update the date heading at the top of the document.
$d.Content.Find.Execute("2023-08-28 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-29 Tuesday", 2)

# Update the multiplication answers in the table. Each populated row of the
# 5-column table is targeted directly by (row, column) so that values which
# happen to collide with each other across cells (e.g. "88x27=2376" is both
# an old value in one cell and a new value in another) cannot cross-match.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text  = "96×74=7104"
$tbl.Cell(1, 2).Range.Text  = "88×27=2376"
$tbl.Cell(1, 3).Range.Text  = "46×18=828"
$tbl.Cell(1, 4).Range.Text  = "91×77=7007"
$tbl.Cell(1, 5).Range.Text  = "59×84=4956"

$tbl.Cell(5, 1).Range.Text  = "95×71=6745"
$tbl.Cell(5, 2).Range.Text  = "60×59=3540"
$tbl.Cell(5, 3).Range.Text  = "76×58=4408"
$tbl.Cell(5, 4).Range.Text  = "12×43=516"
$tbl.Cell(5, 5).Range.Text  = "19×97=1843"

$tbl.Cell(10, 1).Range.Text = "21×19=399"
$tbl.Cell(10, 2).Range.Text = "95×51=4845"
$tbl.Cell(10, 3).Range.Text = "84×14=1176"
$tbl.Cell(10, 4).Range.Text = "27×58=1566"
$tbl.Cell(10, 5).Range.Text = "76×82=6232"

$tbl.Cell(15, 1).Range.Text = "92×96=8832"
$tbl.Cell(15, 2).Range.Text = "18×13=234"
$tbl.Cell(15, 3).Range.Text = "83×85=7055"
$tbl.Cell(15, 4).Range.Text = "34×50=1700"
$tbl.Cell(15, 5).Range.Text = "16×66=1056"

$tbl.Cell(20, 1).Range.Text = "75×29=2175"
$tbl.Cell(20, 2).Range.Text = "71×97=6887"
$tbl.Cell(20, 3).Range.Text = "14×87=1218"
$tbl.Cell(20, 4).Range.Text = "50×50=2500"
$tbl.Cell(20, 5).Range.Text = "42×66=2772"
